# Scheduled runner update: refresh market price / profit figures across
# the per-job "Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# For each affected row we overwrite the cached price/profit columns
# (H..N) with the freshly pulled values. Where a column didn't previously
# contain a value (and the new data has one) the cell simply gets created.
# Where a column previously had a value but the new pull has none, the
# cell content is explicitly cleared so the cell is removed entirely.

$wb = $excel.ActiveWorkbook

$updates = @(
    # ---------------------------------------------------------- ALC ----
    @{ Sheet = "ALC"; Row = 9;  Col = "H"; Val = 129.94737 },
    @{ Sheet = "ALC"; Row = 9;  Col = "I"; Val = 139.15384 },
    @{ Sheet = "ALC"; Row = 9;  Col = "K"; Val = 139.15384 },
    @{ Sheet = "ALC"; Row = 9;  Col = "M"; Val = 29.84616 },

    @{ Sheet = "ALC"; Row = 17; Col = "H"; Val = 296.51514 },
    @{ Sheet = "ALC"; Row = 17; Col = "I"; Val = 390.66666 },
    @{ Sheet = "ALC"; Row = 17; Col = "J"; Val = 287.1 },
    @{ Sheet = "ALC"; Row = 17; Col = "K"; Val = 1171.99998 },
    @{ Sheet = "ALC"; Row = 17; Col = "L"; Val = 861.3000000000001 },
    @{ Sheet = "ALC"; Row = 17; Col = "M"; Val = -1003.99998 },
    @{ Sheet = "ALC"; Row = 17; Col = "N"; Val = -1197.3 },

    @{ Sheet = "ALC"; Row = 40; Col = "H"; Val = 2765.6875 },
    @{ Sheet = "ALC"; Row = 40; Col = "I"; Val = 5400 },
    @{ Sheet = "ALC"; Row = 40; Col = "J"; Val = 1568.2727 },
    @{ Sheet = "ALC"; Row = 40; Col = "K"; Val = 5400 },
    @{ Sheet = "ALC"; Row = 40; Col = "L"; Val = 1568.2727 },
    @{ Sheet = "ALC"; Row = 40; Col = "M"; Val = -5225 },
    @{ Sheet = "ALC"; Row = 40; Col = "N"; Val = -1918.2727 },

    @{ Sheet = "ALC"; Row = 48; Col = "H"; Val = 4500 },
    @{ Sheet = "ALC"; Row = 48; Col = "I"; Val = 4000 },
    @{ Sheet = "ALC"; Row = 48; Col = "J"; Val = 5000 },
    @{ Sheet = "ALC"; Row = 48; Col = "K"; Val = 12000 },
    @{ Sheet = "ALC"; Row = 48; Col = "L"; Val = 15000 },
    @{ Sheet = "ALC"; Row = 48; Col = "M"; Val = -11708 },
    @{ Sheet = "ALC"; Row = 48; Col = "N"; Val = -15584 },

    @{ Sheet = "ALC"; Row = 51; Col = "H"; Val = 2526.5454 },
    @{ Sheet = "ALC"; Row = 51; Col = "I"; Val = 2450 },
    @{ Sheet = "ALC"; Row = 51; Col = "J"; Val = 2570.2856 },
    @{ Sheet = "ALC"; Row = 51; Col = "K"; Val = 2450 },
    @{ Sheet = "ALC"; Row = 51; Col = "L"; Val = 2570.2856 },
    @{ Sheet = "ALC"; Row = 51; Col = "M"; Val = -1966 },
    @{ Sheet = "ALC"; Row = 51; Col = "N"; Val = -3538.2856 },

    @{ Sheet = "ALC"; Row = 56; Col = "H"; Val = 4500 },
    @{ Sheet = "ALC"; Row = 56; Col = "I"; Val = 4000 },
    @{ Sheet = "ALC"; Row = 56; Col = "J"; Val = 5000 },
    @{ Sheet = "ALC"; Row = 56; Col = "K"; Val = 12000 },
    @{ Sheet = "ALC"; Row = 56; Col = "L"; Val = 15000 },
    @{ Sheet = "ALC"; Row = 56; Col = "M"; Val = -11466 },
    @{ Sheet = "ALC"; Row = 56; Col = "N"; Val = -16068 },

    @{ Sheet = "ALC"; Row = 94; Col = "H"; Val = 2451.25 },
    @{ Sheet = "ALC"; Row = 94; Col = "I"; Val = 2451.25 },
    @{ Sheet = "ALC"; Row = 94; Col = "K"; Val = 2451.25 },
    @{ Sheet = "ALC"; Row = 94; Col = "M"; Val = -2000.25 },

    # ---------------------------------------------------------- ARM ----
    @{ Sheet = "ARM"; Row = 45; Col = "H"; Val = 1925.7742 },
    @{ Sheet = "ARM"; Row = 45; Col = "I"; Val = 1761.5 },
    @{ Sheet = "ARM"; Row = 45; Col = "J"; Val = 2780 },
    @{ Sheet = "ARM"; Row = 45; Col = "K"; Val = 1761.5 },
    @{ Sheet = "ARM"; Row = 45; Col = "L"; Val = 2780 },
    @{ Sheet = "ARM"; Row = 45; Col = "M"; Val = -1384.5 },
    @{ Sheet = "ARM"; Row = 45; Col = "N"; Val = -3534 },

    @{ Sheet = "ARM"; Row = 63; Col = "H"; Val = 3767.4211 },
    @{ Sheet = "ARM"; Row = 63; Col = "I"; Val = 2631.3333 },
    @{ Sheet = "ARM"; Row = 63; Col = "J"; Val = 4789.9 },
    @{ Sheet = "ARM"; Row = 63; Col = "K"; Val = 2631.3333 },
    @{ Sheet = "ARM"; Row = 63; Col = "L"; Val = 4789.9 },
    @{ Sheet = "ARM"; Row = 63; Col = "M"; Val = -1945.3333 },
    @{ Sheet = "ARM"; Row = 63; Col = "N"; Val = -6161.9 },

    @{ Sheet = "ARM"; Row = 66; Col = "H"; Val = 3767.4211 },
    @{ Sheet = "ARM"; Row = 66; Col = "I"; Val = 2631.3333 },
    @{ Sheet = "ARM"; Row = 66; Col = "J"; Val = 4789.9 },
    @{ Sheet = "ARM"; Row = 66; Col = "K"; Val = 13156.6665 },
    @{ Sheet = "ARM"; Row = 66; Col = "L"; Val = 23949.5 },
    @{ Sheet = "ARM"; Row = 66; Col = "M"; Val = -9724.666499999999 },
    @{ Sheet = "ARM"; Row = 66; Col = "N"; Val = -30813.5 },

    # ---------------------------------------------------------- BSM ----
    @{ Sheet = "BSM"; Row = 19; Col = "H"; Val = 0 },
    @{ Sheet = "BSM"; Row = 19; Col = "J"; Val = 0 },
    @{ Sheet = "BSM"; Row = 19; Col = "L"; Val = 0 },
    @{ Sheet = "BSM"; Row = 19; Col = "N"; Val = $null },

    @{ Sheet = "BSM"; Row = 35; Col = "H"; Val = 49900 },
    @{ Sheet = "BSM"; Row = 35; Col = "J"; Val = 49900 },
    @{ Sheet = "BSM"; Row = 35; Col = "L"; Val = 49900 },
    @{ Sheet = "BSM"; Row = 35; Col = "N"; Val = -50520 },

    @{ Sheet = "BSM"; Row = 68; Col = "H"; Val = 40295 },
    @{ Sheet = "BSM"; Row = 68; Col = "J"; Val = 40295 },
    @{ Sheet = "BSM"; Row = 68; Col = "L"; Val = 40295 },
    @{ Sheet = "BSM"; Row = 68; Col = "N"; Val = -41917 },

    @{ Sheet = "BSM"; Row = 71; Col = "H"; Val = 40295 },
    @{ Sheet = "BSM"; Row = 71; Col = "J"; Val = 40295 },
    @{ Sheet = "BSM"; Row = 71; Col = "L"; Val = 120885 },
    @{ Sheet = "BSM"; Row = 71; Col = "N"; Val = -128997 },

    @{ Sheet = "BSM"; Row = 82; Col = "H"; Val = 18546.54 },

    @{ Sheet = "BSM"; Row = 85; Col = "H"; Val = 18546.54 },

    @{ Sheet = "BSM"; Row = 100; Col = "H"; Val = 30000 },
    @{ Sheet = "BSM"; Row = 100; Col = "J"; Val = 30000 },
    @{ Sheet = "BSM"; Row = 100; Col = "L"; Val = 30000 },
    @{ Sheet = "BSM"; Row = 100; Col = "N"; Val = -32164 },

    @{ Sheet = "BSM"; Row = 112; Col = "H"; Val = 50000 },
    @{ Sheet = "BSM"; Row = 112; Col = "J"; Val = 50000 },
    @{ Sheet = "BSM"; Row = 112; Col = "L"; Val = 50000 },
    @{ Sheet = "BSM"; Row = 112; Col = "N"; Val = -52954 },

    @{ Sheet = "BSM"; Row = 118; Col = "H"; Val = 69550 },
    @{ Sheet = "BSM"; Row = 118; Col = "J"; Val = 69550 },
    @{ Sheet = "BSM"; Row = 118; Col = "L"; Val = 69550 },
    @{ Sheet = "BSM"; Row = 118; Col = "N"; Val = -72864 },

    # ---------------------------------------------------------- CRP ----
    @{ Sheet = "CRP"; Row = 2;  Col = "H"; Val = 0 },
    @{ Sheet = "CRP"; Row = 2;  Col = "I"; Val = 0 },
    @{ Sheet = "CRP"; Row = 2;  Col = "K"; Val = 0 },
    @{ Sheet = "CRP"; Row = 2;  Col = "M"; Val = $null },

    @{ Sheet = "CRP"; Row = 7;  Col = "H"; Val = 82.111115 },
    @{ Sheet = "CRP"; Row = 7;  Col = "J"; Val = 86.333336 },
    @{ Sheet = "CRP"; Row = 7;  Col = "L"; Val = 86.333336 },
    @{ Sheet = "CRP"; Row = 7;  Col = "N"; Val = -312.333336 },

    @{ Sheet = "CRP"; Row = 52; Col = "H"; Val = 80780 },
    @{ Sheet = "CRP"; Row = 52; Col = "J"; Val = 80780 },
    @{ Sheet = "CRP"; Row = 52; Col = "L"; Val = 80780 },
    @{ Sheet = "CRP"; Row = 52; Col = "N"; Val = -81368 },

    @{ Sheet = "CRP"; Row = 58; Col = "H"; Val = 3498983.5 },
    @{ Sheet = "CRP"; Row = 58; Col = "I"; Val = 6995278 },
    @{ Sheet = "CRP"; Row = 58; Col = "J"; Val = 2688.7693 },
    @{ Sheet = "CRP"; Row = 58; Col = "K"; Val = 6995278 },
    @{ Sheet = "CRP"; Row = 58; Col = "L"; Val = 2688.7693 },
    @{ Sheet = "CRP"; Row = 58; Col = "M"; Val = -6995075 },
    @{ Sheet = "CRP"; Row = 58; Col = "N"; Val = -3094.7693 },

    @{ Sheet = "CRP"; Row = 132; Col = "H"; Val = 2621.75 },
    @{ Sheet = "CRP"; Row = 132; Col = "I"; Val = 1868.0465 },
    @{ Sheet = "CRP"; Row = 132; Col = "J"; Val = 4528.1763 },
    @{ Sheet = "CRP"; Row = 132; Col = "K"; Val = 5604.139499999999 },
    @{ Sheet = "CRP"; Row = 132; Col = "L"; Val = 13584.5289 },
    @{ Sheet = "CRP"; Row = 132; Col = "M"; Val = -3074.139499999999 },
    @{ Sheet = "CRP"; Row = 132; Col = "N"; Val = -18644.5289 },

    @{ Sheet = "CRP"; Row = 134; Col = "H"; Val = 1773.6923 },
    @{ Sheet = "CRP"; Row = 134; Col = "I"; Val = 1505 },
    @{ Sheet = "CRP"; Row = 134; Col = "J"; Val = 2378.25 },
    @{ Sheet = "CRP"; Row = 134; Col = "K"; Val = 4515 },
    @{ Sheet = "CRP"; Row = 134; Col = "L"; Val = 7134.75 },
    @{ Sheet = "CRP"; Row = 134; Col = "M"; Val = -1980 },
    @{ Sheet = "CRP"; Row = 134; Col = "N"; Val = -12204.75 },

    @{ Sheet = "CRP"; Row = 136; Col = "H"; Val = 3498983.5 },
    @{ Sheet = "CRP"; Row = 136; Col = "I"; Val = 6995278 },
    @{ Sheet = "CRP"; Row = 136; Col = "J"; Val = 2688.7693 },
    @{ Sheet = "CRP"; Row = 136; Col = "K"; Val = 20985834 },
    @{ Sheet = "CRP"; Row = 136; Col = "L"; Val = 8066.3079 },
    @{ Sheet = "CRP"; Row = 136; Col = "M"; Val = -20983284 },
    @{ Sheet = "CRP"; Row = 136; Col = "N"; Val = -13166.3079 },

    # ---------------------------------------------------------- CUL ----
    @{ Sheet = "CUL"; Row = 68;  Col = "H"; Val = 3227.5 },
    @{ Sheet = "CUL"; Row = 68;  Col = "I"; Val = 1588.8948 },
    @{ Sheet = "CUL"; Row = 68;  Col = "J"; Val = 5303.067 },
    @{ Sheet = "CUL"; Row = 68;  Col = "K"; Val = 4766.6844 },
    @{ Sheet = "CUL"; Row = 68;  Col = "L"; Val = 15909.201 },
    @{ Sheet = "CUL"; Row = 68;  Col = "M"; Val = -3955.6844 },
    @{ Sheet = "CUL"; Row = 68;  Col = "N"; Val = -17531.201 },

    @{ Sheet = "CUL"; Row = 71;  Col = "H"; Val = 3227.5 },
    @{ Sheet = "CUL"; Row = 71;  Col = "I"; Val = 1588.8948 },
    @{ Sheet = "CUL"; Row = 71;  Col = "J"; Val = 5303.067 },
    @{ Sheet = "CUL"; Row = 71;  Col = "K"; Val = 14300.0532 },
    @{ Sheet = "CUL"; Row = 71;  Col = "L"; Val = 47727.603 },
    @{ Sheet = "CUL"; Row = 71;  Col = "M"; Val = -10244.0532 },
    @{ Sheet = "CUL"; Row = 71;  Col = "N"; Val = -55839.603 },

    @{ Sheet = "CUL"; Row = 107; Col = "H"; Val = 735.8555 },
    @{ Sheet = "CUL"; Row = 107; Col = "I"; Val = 341.42374 },
    @{ Sheet = "CUL"; Row = 107; Col = "J"; Val = 1486.5483 },
    @{ Sheet = "CUL"; Row = 107; Col = "K"; Val = 1024.27122 },
    @{ Sheet = "CUL"; Row = 107; Col = "L"; Val = 4459.644899999999 },
    @{ Sheet = "CUL"; Row = 107; Col = "M"; Val = 895.7287799999999 },
    @{ Sheet = "CUL"; Row = 107; Col = "N"; Val = -8299.644899999999 },

    # ---------------------------------------------------------- GSM ----
    @{ Sheet = "GSM"; Row = 133; Col = "H"; Val = 68856 },
    @{ Sheet = "GSM"; Row = 133; Col = "J"; Val = 68856 },
    @{ Sheet = "GSM"; Row = 133; Col = "L"; Val = 68856 },
    @{ Sheet = "GSM"; Row = 133; Col = "N"; Val = -78976 },

    # ---------------------------------------------------------- LTW ----
    @{ Sheet = "LTW"; Row = 46;  Col = "H"; Val = 899.05554 },
    @{ Sheet = "LTW"; Row = 46;  Col = "I"; Val = 690.0833 },
    @{ Sheet = "LTW"; Row = 46;  Col = "J"; Val = 1317 },
    @{ Sheet = "LTW"; Row = 46;  Col = "K"; Val = 690.0833 },
    @{ Sheet = "LTW"; Row = 46;  Col = "L"; Val = 1317 },
    @{ Sheet = "LTW"; Row = 46;  Col = "M"; Val = -502.0833 },
    @{ Sheet = "LTW"; Row = 46;  Col = "N"; Val = -1693 },

    @{ Sheet = "LTW"; Row = 133; Col = "H"; Val = 47822.777 },
    @{ Sheet = "LTW"; Row = 133; Col = "J"; Val = 47822.777 },
    @{ Sheet = "LTW"; Row = 133; Col = "L"; Val = 47822.777 },
    @{ Sheet = "LTW"; Row = 133; Col = "N"; Val = -52882.777 },

    # ---------------------------------------------------------- WVR ----
    @{ Sheet = "WVR"; Row = 136; Col = "H"; Val = 5170.431 },
    @{ Sheet = "WVR"; Row = 136; Col = "I"; Val = 3105.3225 },
    @{ Sheet = "WVR"; Row = 136; Col = "K"; Val = 9315.967500000001 },
    @{ Sheet = "WVR"; Row = 136; Col = "M"; Val = -6765.967500000001 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $cell = $ws.Range($u.Col + $u.Row)
    if ($null -eq $u.Val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $u.Val
    }
}
